# Apply edits described by the diff:
#  - Add new shared strings / task rows about Footer, Login, Register, etc.
#  - Mark several existing tasks as "Done" (were "In progress")
#  - Fill in owner/deadline/status for rows 16-35
#  - Add "checkpoint-style" notes in column J for several rows
#  - Update the active selection on the sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function Copy-Style($srcAddr, $dstAddr) {
    $ws.Range($srcAddr).Copy() | Out-Null
    $ws.Range($dstAddr).PasteSpecial(-4122) | Out-Null
}

# --- Column style reference cells (already formatted in the sheet) ---
# B column (task description): border + left align + wrap  -> style like B16
# C column (status):           border + center align       -> style like C2
# D column (owner):            border + center align       -> style like D2
# E column (deadline, date):   border + center + date fmt  -> style like E3
# J column (notes, no border): center align                -> style like G2

# --- Two tasks moved from "In progress" to "Done" ---
$ws.Range("C7").Value = "Done"
$ws.Range("C8").Value = "Done"

# --- Fill statuses for rows 16, 17, 19, 21, 23 (owner/deadline already present) ---
foreach ($r in 16,17,19,21,23) {
    $addr = "C$r"
    $ws.Range($addr).Value = "Done"
    Copy-Style "C2" $addr
}

# --- Rows 25 and 26: change status from "In progress" to "Done" ---
$ws.Range("C25").Value = "Done"
$ws.Range("C26").Value = "Done"

# --- Notes in column J for rows 26-31 (entered before the new B-column tasks) ---
$notes = @{
    26 = "Footer"
    27 = "Wishlist subheader"
    28 = "BuyBooks subheader"
    29 = "Login"
    30 = "Register"
    31 = "SuppoertContacts"
}
foreach ($r in 26,27,28,29,30,31) {
    $addr = "J$r"
    $ws.Range($addr).Value = $notes[$r]
    Copy-Style "G2" $addr
}

# --- New tasks for rows 27-31 (column B) ---
$taskNames = @{
    27 = "REACT: Footer"
    28 = "REACT: Wishlist subheader"
    29 = "REACT: BuyBooks subheader"
    30 = "REACT: Login"
    31 = "REACT: Register"
}
$owners = @{
    27 = "Laima"
    28 = "Daniel"
    29 = "Daniel"
    30 = "Aurimas"
    31 = "Aurimas"
}
foreach ($r in 27,28,29,30,31) {
    $ws.Range("B$r").Value = $taskNames[$r]
    Copy-Style "B16" "B$r"

    $ws.Range("D$r").Value = $owners[$r]
    Copy-Style "D2" "D$r"

    # 43118 = 2018-01-18 (Excel 1900 date system serial number)
    $ws.Range("E$r").Value = 43118
    Copy-Style "E3" "E$r"
}

# --- Row 32: new task, already "In progress" ---
$ws.Range("B32").Value = "REACT: Support Contacts"
Copy-Style "B16" "B32"

$ws.Range("C32").Value = "In progress"
Copy-Style "C2" "C32"

$ws.Range("D32").Value = "Birute"
Copy-Style "D2" "D32"

$ws.Range("E32").Value = 43118
Copy-Style "E3" "E32"

# --- Rows 33-35: new tasks (no status yet) ---
$moreTasks = @(
    @{ Row = 33; Task = "REACT: Knygos vienos erdve";                                  Owner = "Jurgis" },
    @{ Row = 34; Task = "REACT: Sulipdyti i viena erdve, kur keicias viduriai";         Owner = "Jurgis" },
    @{ Row = 35; Task = "Knygos: Kaip turi atrodyti paduodama info su knygomis";        Owner = "Aurimas" }
)

foreach ($t in $moreTasks) {
    $r = $t.Row
    $ws.Range("B$r").Value = $t.Task
    Copy-Style "B16" "B$r"

    $ws.Range("D$r").Value = $t.Owner
    Copy-Style "D2" "D$r"

    $ws.Range("E$r").Value = 43118
    Copy-Style "E3" "E$r"
}

# --- Row 20: add a CSS note in column J (added last) ---
$ws.Range("J20").Value = "clear: both;"
Copy-Style "G2" "J20"

# --- Update sheet view (scroll position / selection) ---
$ws.Range("H22").Select() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 19
$win.ScrollColumn = 1

Write-Host "Edits applied"
